$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "20.504.90"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").Value = "1.469.41"
$ws.Range("E3").Value = "  -0.19%  "
$ws.Range("D4").Value = "'1.010"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'0.9759"
$ws.Range("E5").Value = "  -1.54%  "
$ws.Range("D6").Value = "'278.59"
$ws.Range("E6").Value = "  -0.66%  "
$ws.Range("D7").Value = "'0.3656"
$ws.Range("E7").Value = "  -2.19%  "
$ws.Range("D8").Value = "'0.3053"
$ws.Range("E8").Value = "  -5.07%  "
$ws.Range("D9").Value = "'39.48"
$ws.Range("E9").Value = "  -5.37%  "
$ws.Range("D10").Value = "'1.057"
$ws.Range("E10").Value = "  -1.67%  "
$ws.Range("D11").Value = "'0.06631"
$ws.Range("E11").Value = "  -1.94%  "
$ws.Range("D12").Value = "'1.004"
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("D13").Value = "'5.473"
$ws.Range("E13").Value = "  -3.56%  "
$ws.Range("D14").Value = "'17.99"
$ws.Range("E14").Value = "  -3.32%  "
$ws.Range("D15").Value = "'6.189"
$ws.Range("E15").Value = "  -2.55%  "
$ws.Range("B16").Value = "Dai"
$ws.Range("C16").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D16").Value = "'0.9761"
$ws.Range("E16").Value = "  -1.48%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.00001027"
$ws.Range("E17").Value = "  -1.69%  "
$ws.Range("D18").Value = "1.470.25"
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("D19").Value = "'0.05924"
$ws.Range("E19").Value = "  +1.63%  "
$ws.Range("D20").Value = "'69.27"
$ws.Range("E20").Value = "  -5.57%  "
$ws.Range("D21").Value = "'5.446"
$ws.Range("E21").Value = "  -5.32%  "
$ws.Range("D22").Value = "'14.45"
$ws.Range("E22").Value = "  -3.93%  "
$ws.Range("D23").Value = "'10.99"
$ws.Range("E23").Value = "  -2.82%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "'2.244"
$ws.Range("E24").Value = "  -2.95%  "
$ws.Range("B25").Value = "WrappedBTC"
$ws.Range("C25").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D25").Value = "20.541.46"
$ws.Range("E25").Value = "  -0.57%  "
$ws.Range("D26").Value = "'140.53"
$ws.Range("E26").Value = "  +1.24%  "
$ws.Range("D27").Value = "'2.113"
$ws.Range("E27").Value = "  -10.28%  "
$ws.Range("D28").Value = "'17.19"
$ws.Range("E28").Value = "  -3.53%  "
$ws.Range("D29").Value = "1.627.63"
$ws.Range("E29").Value = "  -0.88%  "
$ws.Range("D30").Value = "'114.16"
$ws.Range("E30").Value = "  -0.59%  "
$ws.Range("D31").Value = "'3.949"
$ws.Range("E31").Value = "  -0.80%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'4.945"
$ws.Range("E32").Value = "  -8.86%  "
$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").Value = "'0.07994"
$ws.Range("E33").Value = "  +1.02%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'0.8049"
$ws.Range("E34").Value = "  -5.91%  "
$ws.Range("D35").Value = "'1.544"
$ws.Range("E35").Value = "  -3.22%  "
$ws.Range("D36").Value = "'1.217"
$ws.Range("E36").Value = "  +6.93%  "
$ws.Range("D37").Value = "'0.05794"
$ws.Range("E37").Value = "  -4.33%  "
$ws.Range("D38").Value = "'4.684"
$ws.Range("E38").Value = "  -6.12%  "
$ws.Range("D39").Value = "'0.9755"
$ws.Range("E39").Value = "  -1.80%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.02035"
$ws.Range("E40").Value = "  -2.66%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'7.623"
$ws.Range("E41").Value = "  -3.07%  "
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").Value = "'10.34"
$ws.Range("E42").Value = "  -4.94%  "
$ws.Range("D43").Value = "'0.1875"
$ws.Range("E43").Value = "  -3.00%  "
$ws.Range("D44").Value = "'0.5273"
$ws.Range("E44").Value = "  -3.83%  "
$ws.Range("D45").Value = "'3.513"
$ws.Range("E45").Value = "  -2.57%  "
$ws.Range("D46").Value = "'12.13"
$ws.Range("E46").Value = "  -4.03%  "
$ws.Range("D47").Value = "'118.84"
$ws.Range("E47").Value = "  -1.81%  "
$ws.Range("D48").Value = "'0.5173"
$ws.Range("E48").Value = "  -4.74%  "
$ws.Range("D49").Value = "'1.787"
$ws.Range("E49").Value = "  -3.91%  "
$ws.Range("D50").Value = "'0.06449"
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("D51").Value = "'0.9979"
$ws.Range("E51").Value = "  -0.35%  "
